$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-5: sending cluster changed from FAPs to ECs, values refreshed with new TPM data
# Rows 6-9: new rows added for FAPs as sending cluster (previously only ECs target existed)

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Fgl1"
$ws.Cells.Item(2, 3).Value = "Egfr"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = [double]"1"
$ws.Cells.Item(2, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(2, 7).Value = [double]"0.004391333333333333"
$ws.Cells.Item(2, 8).Value = [double]"0.013174"
$ws.Cells.Item(2, 9).Value = [double]"0.04203587120571539"
$ws.Cells.Item(2, 10).Value = [double]"0.04203587120571539"
$ws.Cells.Item(2, 11).Value = [double]"3"
$ws.Cells.Item(2, 12).Value = [double]"1"
$ws.Cells.Item(2, 13).Value = [double]"1.370876333333333"
$ws.Cells.Item(2, 14).Value = [double]"4.112629"
$ws.Cells.Item(2, 15).Value = [double]"0.01103063309339269"
$ws.Cells.Item(2, 16).Value = [double]"0.01103063309339269"
$ws.Cells.Item(2, 17).Value = [double]"0.006019974938444445"
$ws.Cells.Item(2, 18).Value = [double]"0.054179774446"
$ws.Cells.Item(2, 19).Value = [double]"0.0004636822720313573"
$ws.Cells.Item(2, 20).Value = [double]"0.0004636822720313572"

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Fgl1"
$ws.Cells.Item(3, 3).Value = "Egfr"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = [double]"1"
$ws.Cells.Item(3, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(3, 7).Value = [double]"0.004391333333333333"
$ws.Cells.Item(3, 8).Value = [double]"0.013174"
$ws.Cells.Item(3, 9).Value = [double]"0.04203587120571539"
$ws.Cells.Item(3, 10).Value = [double]"0.04203587120571539"
$ws.Cells.Item(3, 11).Value = [double]"3"
$ws.Cells.Item(3, 12).Value = [double]"1"
$ws.Cells.Item(3, 13).Value = [double]"92.91372433333333"
$ws.Cells.Item(3, 14).Value = [double]"278.741173"
$ws.Cells.Item(3, 15).Value = [double]"0.7476219244149905"
$ws.Cells.Item(3, 16).Value = [double]"0.7476219244149904"
$ws.Cells.Item(3, 17).Value = [double]"0.4080151347891111"
$ws.Cells.Item(3, 18).Value = [double]"3.672136213102"
$ws.Cells.Item(3, 19).Value = [double]"0.03142693892527763"
$ws.Cells.Item(3, 20).Value = [double]"0.03142693892527763"

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Fgl1"
$ws.Cells.Item(4, 3).Value = "Egfr"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = [double]"1"
$ws.Cells.Item(4, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(4, 7).Value = [double]"0.004391333333333333"
$ws.Cells.Item(4, 8).Value = [double]"0.013174"
$ws.Cells.Item(4, 9).Value = [double]"0.04203587120571539"
$ws.Cells.Item(4, 10).Value = [double]"0.04203587120571539"
$ws.Cells.Item(4, 11).Value = [double]"3"
$ws.Cells.Item(4, 12).Value = [double]"1"
$ws.Cells.Item(4, 13).Value = [double]"29.718484"
$ws.Cells.Item(4, 14).Value = [double]"89.155452"
$ws.Cells.Item(4, 15).Value = [double]"0.2391271080585153"
$ws.Cells.Item(4, 16).Value = [double]"0.2391271080585153"
$ws.Cells.Item(4, 17).Value = [double]"0.1305037694053333"
$ws.Cells.Item(4, 18).Value = [double]"1.174533924648"
$ws.Cells.Item(4, 19).Value = [double]"0.01005191631614294"
$ws.Cells.Item(4, 20).Value = [double]"0.01005191631614294"

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Fgl1"
$ws.Cells.Item(5, 3).Value = "Egfr"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = [double]"1"
$ws.Cells.Item(5, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(5, 7).Value = [double]"0.004391333333333333"
$ws.Cells.Item(5, 8).Value = [double]"0.013174"
$ws.Cells.Item(5, 9).Value = [double]"0.04203587120571539"
$ws.Cells.Item(5, 10).Value = [double]"0.04203587120571539"
$ws.Cells.Item(5, 11).Value = [double]"3"
$ws.Cells.Item(5, 12).Value = [double]"1"
$ws.Cells.Item(5, 13).Value = [double]"0.275941"
$ws.Cells.Item(5, 14).Value = [double]"0.827823"
$ws.Cells.Item(5, 15).Value = [double]"0.002220334433101459"
$ws.Cells.Item(5, 16).Value = [double]"0.002220334433101458"
$ws.Cells.Item(5, 17).Value = [double]"0.001211748911333333"
$ws.Cells.Item(5, 18).Value = [double]"0.010905740202"
$ws.Cells.Item(5, 19).Value = [double]"9.333369226346802E-05"
$ws.Cells.Item(5, 20).Value = [double]"9.333369226346801E-05"

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Fgl1"
$ws.Cells.Item(6, 3).Value = "Egfr"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = [double]"1"
$ws.Cells.Item(6, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(6, 7).Value = [double]"0.100075"
$ws.Cells.Item(6, 8).Value = [double]"0.300225"
$ws.Cells.Item(6, 9).Value = [double]"0.9579641287942846"
$ws.Cells.Item(6, 10).Value = [double]"0.9579641287942845"
$ws.Cells.Item(6, 11).Value = [double]"3"
$ws.Cells.Item(6, 12).Value = [double]"1"
$ws.Cells.Item(6, 13).Value = [double]"1.370876333333333"
$ws.Cells.Item(6, 14).Value = [double]"4.112629"
$ws.Cells.Item(6, 15).Value = [double]"0.01103063309339269"
$ws.Cells.Item(6, 16).Value = [double]"0.01103063309339269"
$ws.Cells.Item(6, 17).Value = [double]"0.1371904490583334"
$ws.Cells.Item(6, 18).Value = [double]"1.234714041525"
$ws.Cells.Item(6, 19).Value = [double]"0.01056695082136134"
$ws.Cells.Item(6, 20).Value = [double]"0.01056695082136133"

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Fgl1"
$ws.Cells.Item(7, 3).Value = "Egfr"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = [double]"1"
$ws.Cells.Item(7, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(7, 7).Value = [double]"0.100075"
$ws.Cells.Item(7, 8).Value = [double]"0.300225"
$ws.Cells.Item(7, 9).Value = [double]"0.9579641287942846"
$ws.Cells.Item(7, 10).Value = [double]"0.9579641287942845"
$ws.Cells.Item(7, 11).Value = [double]"3"
$ws.Cells.Item(7, 12).Value = [double]"1"
$ws.Cells.Item(7, 13).Value = [double]"92.91372433333333"
$ws.Cells.Item(7, 14).Value = [double]"278.741173"
$ws.Cells.Item(7, 15).Value = [double]"0.7476219244149905"
$ws.Cells.Item(7, 16).Value = [double]"0.7476219244149904"
$ws.Cells.Item(7, 17).Value = [double]"9.298340962658335"
$ws.Cells.Item(7, 18).Value = [double]"83.68506866392501"
$ws.Cells.Item(7, 19).Value = [double]"0.7161949854897129"
$ws.Cells.Item(7, 20).Value = [double]"0.7161949854897127"

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Fgl1"
$ws.Cells.Item(8, 3).Value = "Egfr"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = [double]"1"
$ws.Cells.Item(8, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(8, 7).Value = [double]"0.100075"
$ws.Cells.Item(8, 8).Value = [double]"0.300225"
$ws.Cells.Item(8, 9).Value = [double]"0.9579641287942846"
$ws.Cells.Item(8, 10).Value = [double]"0.9579641287942845"
$ws.Cells.Item(8, 11).Value = [double]"3"
$ws.Cells.Item(8, 12).Value = [double]"1"
$ws.Cells.Item(8, 13).Value = [double]"29.718484"
$ws.Cells.Item(8, 14).Value = [double]"89.155452"
$ws.Cells.Item(8, 15).Value = [double]"0.2391271080585153"
$ws.Cells.Item(8, 16).Value = [double]"0.2391271080585153"
$ws.Cells.Item(8, 17).Value = [double]"2.9740772863"
$ws.Cells.Item(8, 18).Value = [double]"26.7666955767"
$ws.Cells.Item(8, 19).Value = [double]"0.2290751917423724"
$ws.Cells.Item(8, 20).Value = [double]"0.2290751917423723"

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Fgl1"
$ws.Cells.Item(9, 3).Value = "Egfr"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = [double]"1"
$ws.Cells.Item(9, 6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(9, 7).Value = [double]"0.100075"
$ws.Cells.Item(9, 8).Value = [double]"0.300225"
$ws.Cells.Item(9, 9).Value = [double]"0.9579641287942846"
$ws.Cells.Item(9, 10).Value = [double]"0.9579641287942845"
$ws.Cells.Item(9, 11).Value = [double]"3"
$ws.Cells.Item(9, 12).Value = [double]"1"
$ws.Cells.Item(9, 13).Value = [double]"0.275941"
$ws.Cells.Item(9, 14).Value = [double]"0.827823"
$ws.Cells.Item(9, 15).Value = [double]"0.002220334433101459"
$ws.Cells.Item(9, 16).Value = [double]"0.002220334433101458"
$ws.Cells.Item(9, 17).Value = [double]"0.027614795575"
$ws.Cells.Item(9, 18).Value = [double]"0.248533160175"
$ws.Cells.Item(9, 19).Value = [double]"0.002127000740837991"
$ws.Cells.Item(9, 20).Value = [double]"0.00212700074083799"
